# Apply schema header restructuring for 360-giving-schema-titles workbook
$wb = $excel.ActiveWorkbook

# 1. Remove the "Gazeteer" sheet entirely (closes #29)
$wb.Worksheets.Item("Gazeteer").Delete()

# 2. Activity sheet: drop the "Beneficiary Location:Lat-Long" column (X1),
#    shifting later columns left (closes #38 - geography for point locations)
$wsActivity = $wb.Worksheets.Item("Activity")
$wsActivity.Range("X1").EntireColumn.Delete()

# Classification sheet: rebuild header row with new schema titles
$ws = $wb.Worksheets.Item("Classification")
$ws.Range("A1").Value = "ocid"
$ws.Range("B1").Value = "Activity/id:fundingType"
$ws.Range("C1").Value = "Activity/id:classifications"
$ws.Range("D1").Value = "Vocabulary"
$ws.Range("E1").Value = "Code"
$ws.Range("F1").Value = "Title"
$ws.Range("G1").Value = "Description"
$ws.Range("H1").Value = "URL"
$ws.Range("I1").Value = "Last modified"

# Documents sheet: rebuild header row with new schema titles
$ws = $wb.Worksheets.Item("Documents")
$ws.Range("A1").Value = "ocid"
$ws.Range("B1").Value = "Activity/id:relatedDocument"
$ws.Range("C1").Value = "Identifier"
$ws.Range("D1").Value = "Title"
$ws.Range("E1").Value = "Web Address"
$ws.Range("F1").Value = "Description"
$ws.Range("G1").Value = "Document Type"
$ws.Range("H1").Value = "Last modified"

# Event sheet: rebuild header row with new schema titles
$ws = $wb.Worksheets.Item("Event")
$ws.Range("A1").Value = "ocid"
$ws.Range("B1").Value = "Activity/id:plannedDates"
$ws.Range("C1").Value = "Activity/id:awardDate"
$ws.Range("D1").Value = "Activity/id:actualDates"
$ws.Range("E1").Value = "Title"
$ws.Range("F1").Value = "Start Date"
$ws.Range("G1").Value = "End Date"
$ws.Range("H1").Value = "Duration (months)"
$ws.Range("I1").Value = "Description"
$ws.Range("J1").Value = "Last modified"

# GrantProgramme sheet: rebuild header row with new schema titles
$ws = $wb.Worksheets.Item("GrantProgramme")
$ws.Range("A1").Value = "ocid"
$ws.Range("B1").Value = "Activity/id:grantProgramme"
$ws.Range("C1").Value = "Code"
$ws.Range("D1").Value = "Title"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "URL"
$ws.Range("G1").Value = "Last modified"

# Location sheet: rebuild header row with new schema titles
$ws = $wb.Worksheets.Item("Location")
$ws.Range("A1").Value = "ocid"
$ws.Range("B1").Value = "Activity/id:location"
$ws.Range("C1").Value = "Activity/recipientOrganization[]/id:location"
$ws.Range("D1").Value = "Activity/id:beneficiaryLocation"
$ws.Range("E1").Value = "Activity/fundingOrganization[]/id:location"
$ws.Range("F1").Value = "Identifier"
$ws.Range("G1").Value = "Name"
$ws.Range("H1").Value = "Country Code"
$ws.Range("I1").Value = "Latitude"
$ws.Range("J1").Value = "Longitude"
$ws.Range("K1").Value = "Description"
$ws.Range("L1").Value = "Geographic Code"
$ws.Range("M1").Value = "Geographic Code Type"
$ws.Range("N1").Value = "Last modified"

# Organization sheet: rebuild header row with new schema titles
$ws = $wb.Worksheets.Item("Organization")
$ws.Range("A1").Value = "ocid"
$ws.Range("B1").Value = "Activity/id:recipientOrganization"
$ws.Range("C1").Value = "Activity/id:fundingOrganization"
$ws.Range("D1").Value = "Identifier"
$ws.Range("E1").Value = "Name"
$ws.Range("F1").Value = "Contact Name"
$ws.Range("G1").Value = "Charity Number"
$ws.Range("H1").Value = "Company Number"
$ws.Range("I1").Value = "Street Address"
$ws.Range("J1").Value = "City"
$ws.Range("K1").Value = "County"
$ws.Range("L1").Value = "Country"
$ws.Range("M1").Value = "Postal Code"
$ws.Range("N1").Value = "Phone Number"
$ws.Range("O1").Value = "Alternate Name"
$ws.Range("P1").Value = "Email"
$ws.Range("Q1").Value = "Description"
$ws.Range("R1").Value = "Organisation Type"
$ws.Range("S1").Value = "Web Address"
$ws.Range("T1").Value = "Last modified"

# Transaction sheet: rebuild header row with new schema titles
$ws = $wb.Worksheets.Item("Transaction")
$ws.Range("A1").Value = "ocid"
$ws.Range("B1").Value = "Activity/id:applicationTransaction"
$ws.Range("C1").Value = "Activity/id:commitmentTransaction"
$ws.Range("D1").Value = "Activity/id:disbursementTransaction"
$ws.Range("E1").Value = "Identifier"
$ws.Range("F1").Value = "Transaction date"
$ws.Range("G1").Value = "Currency"
$ws.Range("H1").Value = "Value"
$ws.Range("I1").Value = "Value date"
$ws.Range("J1").Value = "Description"
$ws.Range("K1").Value = "Provider"
$ws.Range("L1").Value = "Recipient"
$ws.Range("M1").Value = "Last modified"
